$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 12.386226927486634
$ws.Range("C2").Value = 11.625879656142747
$ws.Range("D2").Value = 12.749146638718415
$ws.Range("E2").Value = 12.718423230545389

# Row 3 values
$ws.Range("B3").Value = 12.732355385724601
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 13.855620047232103
$ws.Range("E3").Value = 11.054580937996574

# Update selection to match updated used range
$ws.Range("B1:E3").Select()
